$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.291.01"
$ws.Range("E2").Value = "  +6.66%  "
$ws.Range("D3").Value = "2.483.63"
$ws.Range("E3").Value = "  +4.93%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "489.19"
$ws.Range("E5").Value = "  +7.23%  "
$ws.Range("D6").Value = "146.40"
$ws.Range("E6").Value = "  +13.89%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +7.93%  "
$ws.Range("D9").Value = "2.496.20"
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("D10").Value = "5.78"
$ws.Range("E10").Value = "  +11.63%  "
$ws.Range("D11").Value = "0.0982"
$ws.Range("E11").Value = "  +5.12%  "
$ws.Range("E12").Value = "  +7.40%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "2.903.89"
$ws.Range("E14").Value = "  +4.88%  "
$ws.Range("D15").Value = "56.320.97"
$ws.Range("E15").Value = "  +6.59%  "
$ws.Range("D16").Value = "21.21"
$ws.Range("E16").Value = "  +9.75%  "
$ws.Range("E17").Value = "  +6.24%  "
$ws.Range("D18").Value = "2.492.95"
$ws.Range("E18").Value = "  +5.12%  "
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  +11.07%  "
$ws.Range("D20").Value = "10.12"
$ws.Range("E20").Value = "  +9.62%  "
$ws.Range("D21").Value = "318.75"
$ws.Range("E21").Value = "  +4.84%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  +10.87%  "
$ws.Range("D24").Value = "58.61"
$ws.Range("E24").Value = "  +6.36%  "
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  +8.65%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  +8.76%  "
$ws.Range("D28").Value = "2.588.48"
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("D29").Value = "7.63"
$ws.Range("E29").Value = "  +9.40%  "
$ws.Range("D30").Value = "0.0₃0788"
$ws.Range("E30").Value = "  +10.19%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "149.08"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").Value = "18.23"
$ws.Range("E33").Value = "  +4.90%  "
$ws.Range("E34").Value = "  +6.98%  "
$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  +6.40%  "
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  +9.99%  "
$ws.Range("D37").Value = "3.73"
$ws.Range("E37").Value = "  +7.64%  "
$ws.Range("D38").Value = "0.859"
$ws.Range("E38").Value = "  +9.38%  "
$ws.Range("D39").Value = "34.17"
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("D40").Value = "3.52"
$ws.Range("E40").Value = "  +9.30%  "
$ws.Range("D41").Value = "0.0561"
$ws.Range("E41").Value = "  +8.54%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.611"
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.993"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +9.22%  "
$ws.Range("D45").Value = "4.80"
$ws.Range("E45").Value = "  +16.85%  "
$ws.Range("D46").Value = "0.0920"
$ws.Range("E46").Value = "  +7.76%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "258.46"
$ws.Range("E47").Value = "  +19.61%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0229"
$ws.Range("E48").Value = "  +6.77%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "10.18"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").Value = "1.893.70"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "17.62"
$ws.Range("E51").Value = "  +9.33%  "
